$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1) Capture/propagate existing cell formats onto their new homes
#    BEFORE we overwrite the donor cells' own formatting.
#    - Style "4" (black Calibri font + full thin border) currently
#      lives on C2/C3 (and G2/G3). In the reverted layout this style
#      belongs on the "Client" column (E).
# -----------------------------------------------------------------
$ws.Range("C2").Copy() | Out-Null
$ws.Range("E2").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Copy() | Out-Null
$ws.Range("E3").PasteSpecial(-4122) | Out-Null

# Columns C, D and G currently carry styles "4"/"5"; in the reverted
# layout they should look like the plain bordered style ("1"), which
# already lives on cells such as B2/B3.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("C2").PasteSpecial(-4122) | Out-Null
$ws.Range("D2").PasteSpecial(-4122) | Out-Null
$ws.Range("G2").PasteSpecial(-4122) | Out-Null
$ws.Range("B3").Copy() | Out-Null
$ws.Range("C3").PasteSpecial(-4122) | Out-Null
$ws.Range("D3").PasteSpecial(-4122) | Out-Null
$ws.Range("G3").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# -----------------------------------------------------------------
# 2) Header row (row 1)
# -----------------------------------------------------------------
$ws.Range("A1").Value = "Order Received Data and Time"
$ws.Range("B1").Value = "OrderID"
$ws.Range("C1").Value = "Emp ID-Order Assigned"
$ws.Range("D1").Value = "Assignee_QA"
$ws.Range("E1").Value = "Client"
$ws.Range("F1").Value = "Typist"
$ws.Range("G1").Value = "Typist QC"
$ws.Range("H1").Value = "Product Name"
$ws.Range("I1").Value = "Process"
$ws.Range("J1").Value = "Lob"
$ws.Range("K1").Value = "State"
$ws.Range("L1").Value = "County"
$ws.Range("M1").Value = "Status"

# -----------------------------------------------------------------
# 3) Data row 2
# -----------------------------------------------------------------
$ws.Range("A2").Value = 45436.041666666664
$ws.Range("B2").Value = 121321783
$ws.Range("C2").Value = "SIPL0005"
$ws.Range("D2").Value = "SIPL0004"
$ws.Range("E2").Value = "Qualia"
$ws.Range("F2").Value = "SIPL0102"
$ws.Range("G2").Value = "SIPL0103"
$ws.Range("H2").Value = "Current Owner Search"
$ws.Range("I2").Value = "Search & Typing"
$ws.Range("J2").Value = "Title"
$ws.Range("K2").Value = "FL"
$ws.Range("L2").Value = "Clay"
$ws.Range("M2").Value = "WIP"

# -----------------------------------------------------------------
# 4) Data row 3
# -----------------------------------------------------------------
$ws.Range("A3").Value = 45439.083333333336
$ws.Range("B3").Value = 2193218321
$ws.Range("C3").Value = "SIPL0005"
$ws.Range("D3").Value = "SIPL0004"
$ws.Range("E3").Value = "Qualia"
$ws.Range("F3").Value = "SIPL0102"
$ws.Range("G3").Value = "SIPL0103"
$ws.Range("H3").Value = "Full Search"
$ws.Range("I3").Value = "Search & Typing"
$ws.Range("J3").Value = "Title"
$ws.Range("K3").Value = "FL"
$ws.Range("L3").Value = "Clay"
$ws.Range("M3").Value = "WIP"

# -----------------------------------------------------------------
# 5) Drop the now-unused "Tier" column (N) entirely
# -----------------------------------------------------------------
$ws.Range("N1:N3").Delete(-4159) | Out-Null

# -----------------------------------------------------------------
# 6) Column widths to match the reverted layout
#    (values chosen so the engine's internal pixel-rounding lands as
#    close as possible on the widths recorded in the target file)
# -----------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 35.5
$ws.Columns.Item(5).ColumnWidth = 11.76
$ws.Columns.Item(8).ColumnWidth = 15.14
$ws.Columns.Item(9).ColumnWidth = 15.14
$ws.Columns.Item(10).ColumnWidth = 15.14

# -----------------------------------------------------------------
# 7) Selection matches the reverted workbook state
# -----------------------------------------------------------------
$ws.Range("H12").Select() | Out-Null
